$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "v" markers that filled column G (rows 4-16)
$ws.Range("G4:G16").ClearContents()

# Build out new rows 23..28 (copying the formatting of row 22 down) so
# that row 23 becomes a real data row and row 28 ends up as a blank,
# styled row - mirroring what the workbook looks like after the table
# was extended in the UI.
$ws.Range("D22:G22").Copy()
$ws.Range("D23:G23").Insert(-4121)
$ws.Range("D23:G23").Copy()
$ws.Range("D24:G24").Insert(-4121)
$ws.Range("D24:G24").Copy()
$ws.Range("D25:G25").Insert(-4121)
$ws.Range("D25:G25").Copy()
$ws.Range("D26:G26").Insert(-4121)
$ws.Range("D26:G26").Copy()
$ws.Range("D27:G27").Insert(-4121)
$ws.Range("D27:G27").Copy()
$ws.Range("D28:G28").Insert(-4121)

# Rows 24-27 should stay empty/unused - drop the placeholder cells that
# the repeated insert created for them.
$ws.Range("D24:G27").Clear()

# Fill in the new dictionary entry on row 23.
$ws.Range("D23").Value = "numéro de semaine"
$ws.Range("E23").Value = "week_id"
$ws.Range("F23").Value = "entier"

# Extend the table / autofilter so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D3:G28"))

# Match the cursor position left behind in the workbook.
$ws.Range("G5").Select()
